$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look numeric (e.g. "215.52") but the source
# data stores them as plain text (inline strings). Excel's COM layer
# auto-detects numeric-looking strings and would coerce them to real
# numbers on assignment, so we temporarily force the Text number format
# on the whole Price column before writing, then restore the default
# "Normal" style afterwards so no visible formatting changes remain.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "25.803.71"
$ws.Range("D3").Value = "1.636.80"
$ws.Range("D5").Value = "215.52"
$ws.Range("D6").Value = "0.505"
$ws.Range("D9").Value = "0.0638"
$ws.Range("D10").Value = "19.75"
$ws.Range("D11").Value = "0.0793"
$ws.Range("D13").Value = "1.861.64"
$ws.Range("D14").Value = "1.634.29"
$ws.Range("D17").Value = "63.23"
$ws.Range("D18").Value = "25.821.34"
$ws.Range("D20").Value = "4.47"
$ws.Range("D21").Value = "192.83"
$ws.Range("D22").Value = "9.99"
$ws.Range("D23").Value = "6.37"
$ws.Range("D26").Value = "142.47"
$ws.Range("D28").Value = "6.96"
$ws.Range("D29").Value = "15.55"
$ws.Range("D31").Value = "0.0495"
$ws.Range("D37").Value = "1.132.54"
$ws.Range("D39").Value = "0.545"
$ws.Range("D43").Value = "100.57"
$ws.Range("D44").Value = "0.806"
$ws.Range("D45").Value = "1.771.24"
$ws.Range("D47").Value = "55.32"

$priceRange.Style = "Normal"

# Column E ("Volume(1h)") values are percentages formatted as
# plain text with surrounding padding spaces; they are not numeric
# so they can be assigned directly.
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +3.18%  "
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("E51").Value = "  +2.67%  "
